$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row 1 - publisher text (merged A1:B1)
$ws.Range("A1").Value = "Publisher: The Royal Society of Chemistry"

# Row 2 - add ISSN column headers
$ws.Range("C2").Value = "ISSN(print)"
$ws.Range("D2").Value = "ISSN(online)"
$ws.Range("C2").Style = $ws.Range("A2").Style
$ws.Range("D2").Style = $ws.Range("A2").Style

# Row 41 - Physical Chemistry Chemical Physics ISSN values
$ws.Range("C41").Value = "1463-9076"
$ws.Range("D41").Value = "1463-9084"

# Update selection to match target
$ws.Range("A1:B1").Select()
